$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1288.4286
$ws.Range("I19").Value = 785
$ws.Range("K19").Value = 785
$ws.Range("M19").Value = -610

$ws.Range("H113").Value = 71497.94
$ws.Range("I113").Value = 175904.33
$ws.Range("J113").Value = 14549
$ws.Range("K113").Value = 175904.33
$ws.Range("L113").Value = 14549
$ws.Range("M113").Value = -172650.33
$ws.Range("N113").Value = -21057

$ws.Range("H116").Value = 9374
$ws.Range("I116").Value = 11747.083
$ws.Range("K116").Value = 11747.083
$ws.Range("M116").Value = -8305.083000000001

$ws.Range("H125").Value = 1933
$ws.Range("I125").Value = 2391
$ws.Range("J125").Value = 1475
$ws.Range("K125").Value = 21519
$ws.Range("L125").Value = 13275
$ws.Range("M125").Value = -19059
$ws.Range("N125").Value = -18195

$ws.Range("H135").Value = 850.5833
$ws.Range("I135").Value = 700.6667
$ws.Range("K135").Value = 6306.0003
$ws.Range("M135").Value = -3771.0003

$ws.Range("H137").Value = 6243.3335
$ws.Range("I137").Value = 1399.2
$ws.Range("K137").Value = 4197.6
$ws.Range("M137").Value = -1647.6

$ws.Range("H138").Value = 3004.0308
$ws.Range("I138").Value = 1239.9445
$ws.Range("K138").Value = 3719.8335
$ws.Range("M138").Value = 1420.1665

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 0
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("N37").ClearContents()

$ws.Range("H44").Value = 19832.666

$ws.Range("H74").Value = 35579.145
$ws.Range("I74").Value = 57205.332
$ws.Range("K74").Value = 57205.332
$ws.Range("M74").Value = -56331.332

$ws.Range("H77").Value = 35579.145
$ws.Range("I77").Value = 57205.332
$ws.Range("K77").Value = 286026.66
$ws.Range("M77").Value = -281658.66

$ws.Range("H96").Value = 8750
$ws.Range("J96").Value = 8750
$ws.Range("L96").Value = 8750
$ws.Range("N96").Value = -14242

$ws.Range("H106").Value = 216666.67
$ws.Range("J106").Value = 216666.67
$ws.Range("L106").Value = 216666.67
$ws.Range("N106").Value = -219190.67

$ws.Range("H125").Value = 29374.125
$ws.Range("J125").Value = 29374.125
$ws.Range("L125").Value = 29374.125
$ws.Range("N125").Value = -39214.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2102.7144
$ws.Range("I86").Value = 1948.909
$ws.Range("J86").Value = 2666.6667
$ws.Range("K86").Value = 1948.909
$ws.Range("L86").Value = 2666.6667
$ws.Range("M86").Value = -825.9090000000001
$ws.Range("N86").Value = -4912.6667

$ws.Range("H89").Value = 2102.7144
$ws.Range("I89").Value = 1948.909
$ws.Range("J89").Value = 2666.6667
$ws.Range("K89").Value = 9744.545
$ws.Range("L89").Value = 13333.3335
$ws.Range("M89").Value = -4128.545
$ws.Range("N89").Value = -24565.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 25965.023
$ws.Range("I31").Value = 50682.145
$ws.Range("K31").Value = 50682.145
$ws.Range("M31").Value = -50387.145

$ws.Range("H34").Value = 25965.023
$ws.Range("I34").Value = 50682.145
$ws.Range("K34").Value = 50682.145
$ws.Range("M34").Value = -50480.145

$ws.Range("H37").Value = 29999
$ws.Range("J37").Value = 29999
$ws.Range("L37").Value = 29999
$ws.Range("N37").Value = -30213

$ws.Range("H58").Value = 2587.2222
$ws.Range("I58").Value = 2445.7368
$ws.Range("J58").Value = 3355.2856
$ws.Range("K58").Value = 2445.7368
$ws.Range("L58").Value = 3355.2856
$ws.Range("M58").Value = -2242.7368
$ws.Range("N58").Value = -3761.2856

$ws.Range("H70").Value = 55000
$ws.Range("J70").Value = 55000
$ws.Range("L70").Value = 55000
$ws.Range("N70").Value = -55630

$ws.Range("H73").Value = 55000
$ws.Range("J73").Value = 55000
$ws.Range("L73").Value = 55000
$ws.Range("N73").Value = -57184

$ws.Range("H127").Value = 83179.75
$ws.Range("I127").Value = 79666.664
$ws.Range("K127").Value = 79666.664
$ws.Range("M127").Value = -74706.664

$ws.Range("H132").Value = 113557.53
$ws.Range("I132").Value = 187099.92
$ws.Range("K132").Value = 561299.76
$ws.Range("M132").Value = -558769.76

$ws.Range("H134").Value = 18331.688
$ws.Range("I134").Value = 14676.615
$ws.Range("K134").Value = 44029.845
$ws.Range("M134").Value = -41494.845

$ws.Range("H136").Value = 2587.2222
$ws.Range("I136").Value = 2445.7368
$ws.Range("J136").Value = 3355.2856
$ws.Range("K136").Value = 7337.2104
$ws.Range("L136").Value = 10065.8568
$ws.Range("M136").Value = -4787.2104
$ws.Range("N136").Value = -15165.8568

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 23333.334
$ws.Range("J15").Value = 23333.334
$ws.Range("L15").Value = 23333.334
$ws.Range("N15").Value = -23909.334

$ws.Range("H81").Value = 23333.334
$ws.Range("J81").Value = 23333.334
$ws.Range("L81").Value = 23333.334
$ws.Range("N81").Value = -25329.334

$ws.Range("H84").Value = 23333.334
$ws.Range("J84").Value = 23333.334
$ws.Range("L84").Value = 70000.00199999999
$ws.Range("N84").Value = -79984.00199999999

$ws.Range("H132").Value = 2922.551
$ws.Range("I132").Value = 2858.2
$ws.Range("K132").Value = 8574.599999999999
$ws.Range("M132").Value = -6044.599999999999

$ws.Range("H134").Value = 25757.691
$ws.Range("J134").Value = 25757.691
$ws.Range("L134").Value = 77273.073
$ws.Range("N134").Value = -82343.073

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 28474.75
$ws.Range("I22").Value = 1949.5
$ws.Range("K22").Value = 1949.5
$ws.Range("M22").Value = -1654.5

$ws.Range("H27").Value = 28474.75
$ws.Range("I27").Value = 1949.5
$ws.Range("K27").Value = 1949.5
$ws.Range("M27").Value = -1842.5

$ws.Range("H43").Value = 10666.667
$ws.Range("J43").Value = 10666.667
$ws.Range("L43").Value = 10666.667
$ws.Range("N43").Value = -11052.667

$ws.Range("H46").Value = 1089
$ws.Range("I46").Value = 1130
$ws.Range("K46").Value = 1130
$ws.Range("M46").Value = -942

$ws.Range("H68").Value = 4031.8
$ws.Range("I68").Value = 4616.857
$ws.Range("J68").Value = 2666.6667
$ws.Range("K68").Value = 4616.857
$ws.Range("L68").Value = 2666.6667
$ws.Range("M68").Value = -3867.857
$ws.Range("N68").Value = -4164.6667

$ws.Range("H71").Value = 4031.8
$ws.Range("I71").Value = 4616.857
$ws.Range("J71").Value = 2666.6667
$ws.Range("K71").Value = 23084.285
$ws.Range("L71").Value = 13333.3335
$ws.Range("M71").Value = -19340.285
$ws.Range("N71").Value = -20821.3335

$ws.Range("H93").Value = 2261.625
$ws.Range("J93").Value = 1847.5
$ws.Range("L93").Value = 1847.5
$ws.Range("N93").Value = -4343.5

$ws.Range("H100").Value = 1860
$ws.Range("I100").Value = 1860
$ws.Range("K100").Value = 1860
$ws.Range("M100").Value = -1319

$ws.Range("H135").Value = 207614.47
$ws.Range("J135").Value = 207614.47
$ws.Range("L135").Value = 207614.47
$ws.Range("N135").Value = -217754.47

$ws.Range("H136").Value = 2312.8
$ws.Range("I136").Value = 1966.6552
$ws.Range("J136").Value = 3225.3635
$ws.Range("K136").Value = 5899.9656
$ws.Range("L136").Value = 9676.0905
$ws.Range("M136").Value = -3349.9656
$ws.Range("N136").Value = -14776.0905

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 29000
$ws.Range("J39").Value = 29000
$ws.Range("L39").Value = 29000
$ws.Range("N39").Value = -29826

$ws.Range("H42").Value = 28200
$ws.Range("J42").Value = 28200
$ws.Range("L42").Value = 28200
$ws.Range("N42").Value = -28956

$ws.Range("H43").Value = 13600
$ws.Range("I43").Value = 10000
$ws.Range("J43").Value = 28000
$ws.Range("K43").Value = 10000
$ws.Range("L43").Value = 28000
$ws.Range("M43").Value = -9851
$ws.Range("N43").Value = -28298

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
